$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-05-15", "12:35:28", "-", "Etiquetadora", "-", "-", "-", "12:35:53"),
    @("2024-05-15", "12:37:22", "Palet atascado en la curva", "-", "-", "-", "-", "12:37:29"),
    @("2024-05-15", "12:37:26", "Fallo fijador tapa", "-", "-", "-", "-", "12:37:30"),
    @("2024-05-15", "12:38:22", "-", "-", "-", "-", "Soldadura defectuosa", "12:38:28"),
    @("2024-05-15", "12:38:37", "-", "-", "-", "-", "Marco atascado en parte inferior", $null),
    @("2024-05-15", "12:46:36", "AOI no detecta pieza", "-", "-", "-", "-", "12:46:48")
)

$startRow = 60
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $cell = $ws.Cells.Item($r, $c + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        }
    }
}
